# "Updated Master data as per 16th May Refresh"
#
# Appends 3 new user_detail rows (ids 110033-110035) below the existing
# 32 data rows (sheet currently spans A1:K33 -> A1:K36 afterwards), using
# the same column layout / constants as every prior row:
#   status_code=ACT, lang_code=eng, last_login_method=PWD, is_active=TRUE,
#   cr_by=superadmin, cr_dtimes=now()

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 34-36 by duplicating the formatting of the last existing row
# (33) - this carries over the left-aligned "email" column style (D) and
# the is_active boolean column style (I) without minting any new styles.
$ws.Range("A33:K33").Copy()
$ws.Range("A34:K36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @{ id = 110033; uin = 9317596771; name = "Nikola Tesla"; email = "nikola.tesla@xyz.com"; mobile = 818876434 },
    @{ id = 110034; uin = 9317596772; name = "Graham Bell";  email = "graham.bell@xyz.com";  mobile = 818876435 },
    @{ id = 110035; uin = 9317596773; name = "Albert Miles"; email = "albert.miles@xyz.com"; mobile = 818876436 }
)
$startRow = 34

# Write column-by-column (id, then uin, then name, then email, ...)
# rather than row-by-row, so the new shared-string entries are appended
# in the same order as the source workbook: all 3 names first, then all
# 3 emails.
$r = $startRow
foreach ($row in $newRows) { $ws.Cells.Item($r, 1).Value = $row.id; $r = $r + 1 }

$r = $startRow
foreach ($row in $newRows) { $ws.Cells.Item($r, 2).Value = $row.uin; $r = $r + 1 }

$r = $startRow
foreach ($row in $newRows) { $ws.Cells.Item($r, 3).Value = $row.name; $r = $r + 1 }

$r = $startRow
foreach ($row in $newRows) { $ws.Cells.Item($r, 4).Value = $row.email; $r = $r + 1 }

$r = $startRow
foreach ($row in $newRows) { $ws.Cells.Item($r, 5).Value = $row.mobile; $r = $r + 1 }

$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 6).Value  = "ACT"
    $ws.Cells.Item($r, 7).Value  = "eng"
    $ws.Cells.Item($r, 8).Value  = "PWD"
    $ws.Cells.Item($r, 9).Value  = $true
    $ws.Cells.Item($r, 10).Value = "superadmin"
    $ws.Cells.Item($r, 11).Value = "now()"
    $r = $r + 1
}

# Leave the sheet selection on the row below the newly entered data, as
# Excel would after typing values into the last populated row.
$ws.Range("A37:XFD1048576").Select()
